$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new price records (week of 2023-04-27, serial 45043) are being added to
# the top of the data table. Insert two blank rows right after the header
# block of unchanged rows (row 9) which pushes every existing record down by
# two rows; the new rows then get filled in with the new data.
$ws.Rows("10:11").Insert()

$row10 = @(8, "Terminal La Palmera de La Serena", "Coquimbo", 45043, 4, "Fruta", 100104, "Frutos de pepita", 100104003, "Membrillo", "Champion", "Especial", 10, 320000, 330000, 325000, "`$/bins (450 kilos)", "Región de O'Higgins", 722, 450)
$row11 = @(8, "Terminal La Palmera de La Serena", "Coquimbo", 45043, 4, "Fruta", 100104, "Frutos de pepita", 100104003, "Membrillo", "Champion", "Primera", 20, 290000, 300000, 295000, "`$/bins (450 kilos)", "Región de O'Higgins", 656, 450)

for ($i = 0; $i -lt $row10.Length; $i++) {
    $ws.Cells.Item(10, $i + 1).Value = $row10[$i]
}
for ($i = 0; $i -lt $row11.Length; $i++) {
    $ws.Cells.Item(11, $i + 1).Value = $row11[$i]
}
